# Update the "Datos" sheet in the e-prepago datadriven workbook:
#  - D2 (usuario) test account name bumped from autotest29 -> autotest32
#  - B2 (numeroDocumento) bumped from 333333301 -> 333333304
#  - active cell selection moved from C2 to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

$ws.Range("D2").Value = "autotest32"
$ws.Range("B2").Value = 333333304

$ws.Range("B2").Select()
